$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Valor Total: 96"
$ws.Range("A5:D5").Merge()
$ws.Range("A5").Style = $ws.Range("A1").Style
$ws.Range("B5:D5").Style = $ws.Range("B1:D1").Style
